# Framework creation in Progress
# Applies the TestCaseInfo table edits (renumbered/rearranged test rows,
# new Test8 row content) plus the selection changes on the two sheets.

$wb = $excel.ActiveWorkbook

# --- TestCaseInfo sheet: update the test case table ---
$ws1 = $wb.Worksheets.Item("TestCaseInfo")

# New row 9 (Test8) first, so its new shared string is introduced
# before the "ignore" values that now appear higher up in the sheet.
$ws1.Range("A9").Value = 8
$ws1.Range("B9").Value = "Test8"
$ws1.Range("C9").Value = "frameworkcore.projectA.moduleA.featureA.FeatureBTest"
$ws1.Range("F9").Value = "N"

$ws1.Range("C2").Value = "com.cucumber.Test.TestRunnerA"
$ws1.Range("D2").Value = "ignore"
$ws1.Range("F2").Value = "Y"

$ws1.Range("B3").Value = "Test2"
$ws1.Range("C3").Value = "com.cucumber.Test.TestRunnerB"
$ws1.Range("D3").Value = "ignore"
$ws1.Range("E3").Value = "Firefox"

$ws1.Range("B4").Value = "Test3"
$ws1.Range("C4").Value = "com.api.tests.module1.Module1Test"

$ws1.Range("B5").Value = "Test4"
$ws1.Range("C5").Value = "ModuleC"

$ws1.Range("B6").Value = "Test5"
$ws1.Range("C6").Value = "ModuleD"

$ws1.Range("B7").Value = "Test6"
$ws1.Range("C7").Value = "frameworkcore.projectA.moduleA.featureA.FeatureATest"

$ws1.Range("B8").Value = "Test7"
$ws1.Range("C8").Value = "frameworkcore.projectA.moduleA.featureA.FeatureATest"
$ws1.Range("F8").Value = "N"

# Restore the selected cell on this sheet
$ws1.Activate() | Out-Null
$ws1.Range("F3").Select() | Out-Null

# --- Configuration sheet: just a selection change ---
$ws2 = $wb.Worksheets.Item("Configuration")
$ws2.Activate() | Out-Null
$ws2.Range("A2").Select() | Out-Null
